$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.940.52'
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("D3").Value = '2.592.83'
$ws.Range("E3").Value = '  -0.80%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.22'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.89%  '
$ws.Range("E7").Value = '  -1.12%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.26'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.87%  '
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("E13").Value = '  -1.56%  '
$ws.Range("D14").Value = '2.993.12'
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("E15").Value = '  +0.98%  '
$ws.Range("D16").Value = '2.585.64'
$ws.Range("E16").Value = '  -1.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.912'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.79'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("D19").Value = '46.099.32'
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.49%  '
$ws.Range("B22").Value = 'InternetComputer(DFINITY)'
$ws.Range("C22").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.73'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '292.44'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +13.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.78'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.03'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.56%  '
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("E27").Value = '  +4.34%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.04'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.73'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '38.17'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.59%  '
$ws.Range("E32").Value = '  -2.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.21'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.56'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '154.96'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0833'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.35%  '
$ws.Range("E37").Value = '  -6.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.77'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.87%  '
$ws.Range("E39").Value = '  +2.83%  '
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.62'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.35%  '
$ws.Range("E42").Value = '  +1.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.55'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.16'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.92'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.38%  '
$ws.Range("D46").Value = '2.105.23'
$ws.Range("E46").Value = '  +2.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '97.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.58%  '
$ws.Range("E48").Value = '  +0.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.61'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.08%  '
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '107.88'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.19%  '
